$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new columns I0 (I) and IF (J) with header formatting matching existing headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I0 and IF values for rows 2-88
$ijValues = @{
    2 = @(7, 8)
    3 = @(8, 10)
    4 = @(9, 9)
    5 = @(8, 9)
    6 = @(9, 9)
    7 = @(9, 9)
    8 = @(9, 9)
    9 = @(9, 9)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(10, 10)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(9, 9)
    22 = @(10, 10)
    23 = @(9, 9)
    24 = @(8, 8)
    25 = @(9, 9)
    26 = @(9, 9)
    27 = @(11, 11)
    28 = @(9, 9)
    29 = @(9, 9)
    30 = @(9, 9)
    31 = @(9, 9)
    32 = @(9, 9)
    33 = @(9, 9)
    34 = @(9, 9)
    35 = @(9, 9)
    36 = @(9, 9)
    37 = @(9, 9)
    38 = @(9, 9)
    39 = @(9, 9)
    40 = @(9, 9)
    41 = @(9, 9)
    42 = @(9, 9)
    43 = @(9, 9)
    44 = @(9, 9)
    45 = @(9, 9)
    46 = @(8, 8)
    47 = @(9, 9)
    48 = @(9, 9)
    49 = @(9, 9)
    50 = @(9, 9)
    51 = @(9, 9)
    52 = @(11, 11)
    53 = @(9, 9)
    54 = @(9, 9)
    55 = @(9, 9)
    56 = @(10, 10)
    57 = @(9, 9)
    58 = @(9, 9)
    59 = @(9, 9)
    60 = @(9, 9)
    61 = @(9, 9)
    62 = @(9, 9)
    63 = @(9, 9)
    64 = @(9, 9)
    65 = @(9, 9)
    66 = @(9, 9)
    67 = @(9, 9)
    68 = @(10, 10)
    69 = @(10, 10)
    70 = @(9, 9)
    71 = @(9, 9)
    72 = @(9, 9)
    73 = @(9, 9)
    74 = @(9, 9)
    75 = @(9, 9)
    76 = @(8, 8)
    77 = @(9, 9)
    78 = @(9, 9)
    79 = @(9, 9)
    80 = @(9, 9)
    81 = @(9, 9)
    82 = @(10, 10)
    83 = @(5, 5)
    84 = @(5, 5)
    85 = @(7, 7)
    86 = @(5, 5)
    87 = @(4, 4)
    88 = @(5, 5)
}

foreach ($r in $ijValues.Keys) {
    $vals = $ijValues[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
